# Update the workbook with new TPM-derived values.
# - Row 4's Target cluster label changes from "MuSCs" to "Inflammatory-Mac"
# - Row 5's Target cluster label changes from "Resolving-Mac" to "MuSCs"
# - Several numeric columns (M-T) are refreshed with newly computed TPM values
#   for rows 2-5, and K/L swap between rows 4 and 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Target cluster labels (column D) ---
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("D5").Value = "MuSCs"

# --- Row 2 (ECs) ---
$ws.Range("M2").Value = 0.1579376666666667
$ws.Range("N2").Value = 0.473813
$ws.Range("O2").Value = 0.05467876644486869
$ws.Range("P2").Value = 0.07340983674118848
$ws.Range("Q2").Value = 0.0613432003168889
$ws.Range("R2").Value = 0.552088802852
$ws.Range("S2").Value = 0.05467876644486869
$ws.Range("T2").Value = 0.07340983674118848

# --- Row 3 (FAPs) ---
$ws.Range("O3").Value = 0.1724539210166233
$ws.Range("P3").Value = 0.2315307204300726
$ws.Range("S3").Value = 0.1724539210166233
$ws.Range("T3").Value = 0.2315307204300726

# --- Row 4 (now Inflammatory-Mac) ---
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02135966666666667
$ws.Range("N4").Value = 0.064079
$ws.Range("O4").Value = 0.00739481752299059
$ws.Range("P4").Value = 0.0099280284174107
$ws.Range("Q4").Value = 0.008296123012888889
$ws.Range("R4").Value = 0.07466510711600001
$ws.Range("S4").Value = 0.00739481752299059
$ws.Range("T4").Value = 0.0099280284174107

# --- Row 5 (now MuSCs) ---
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.21104
$ws.Range("N5").Value = 4.422079999999999
$ws.Range("O5").Value = 0.7654724950155174
$ws.Range("P5").Value = 0.6851314144113283
$ws.Range("Q5").Value = 0.8587708840533333
$ws.Range("R5").Value = 5.15262530432
$ws.Range("S5").Value = 0.7654724950155174
$ws.Range("T5").Value = 0.6851314144113283
